$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("G2").Value = 4.3
$ws.Range("H2").Value = 2.26
$ws.Range("I2").Value = 2.76
$ws.Range("J2").Value = 2.82
$ws.Range("P2").Value = 1.57

# Row 3 updates
$ws.Range("G3").Value = 1.75
$ws.Range("K3").Value = 4.2
$ws.Range("Q3").Value = 2.18
